# Auto-generated edit script applying value updates to Ultima Profits market-data sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1191.5883
$ws.Range("I19").Value = 292.3
$ws.Range("J19").Value = 2476.2856
$ws.Range("K19").Value = 292.3
$ws.Range("L19").Value = 2476.2856
$ws.Range("M19").Value = -117.3
$ws.Range("N19").Value = -2826.2856
$ws.Range("H62").Value = 1186.5555
$ws.Range("I62").Value = 1215.8
$ws.Range("J62").Value = 1150
$ws.Range("K62").Value = 1215.8
$ws.Range("L62").Value = 1150
$ws.Range("M62").Value = -591.8
$ws.Range("N62").Value = -2398
$ws.Range("H65").Value = 1186.5555
$ws.Range("I65").Value = 1215.8
$ws.Range("J65").Value = 1150
$ws.Range("K65").Value = 6079
$ws.Range("L65").Value = 5750
$ws.Range("M65").Value = -2959
$ws.Range("N65").Value = -11990
$ws.Range("H116").Value = 2515.3076
$ws.Range("I116").Value = 2322.111
$ws.Range("J116").Value = 2950
$ws.Range("K116").Value = 2322.111
$ws.Range("L116").Value = 2950
$ws.Range("M116").Value = 1119.889
$ws.Range("N116").Value = -9834
$ws.Range("H141").Value = 2081.3635
$ws.Range("I141").Value = 1989.5
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 5968.5
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -788.5
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("H45").Value = 1749524.2
$ws.Range("I45").Value = 2067411.4
$ws.Range("J45").Value = 1145
$ws.Range("K45").Value = 2067411.4
$ws.Range("L45").Value = 1145
$ws.Range("M45").Value = -2067034.4
$ws.Range("N45").Value = -1899
$ws.Range("H74").Value = 2145.0908
$ws.Range("I74").Value = 1932.4445
$ws.Range("J74").Value = 2292.3076
$ws.Range("K74").Value = 1932.4445
$ws.Range("L74").Value = 2292.3076
$ws.Range("M74").Value = -1058.4445
$ws.Range("N74").Value = -4040.3076
$ws.Range("H77").Value = 2145.0908
$ws.Range("I77").Value = 1932.4445
$ws.Range("J77").Value = 2292.3076
$ws.Range("K77").Value = 9662.2225
$ws.Range("L77").Value = 11461.538
$ws.Range("M77").Value = -5294.2225
$ws.Range("N77").Value = -20197.538
$ws.Range("H109").Value = 32722.545
$ws.Range("J109").Value = 32722.545
$ws.Range("L109").Value = 32722.545
$ws.Range("N109").Value = -35496.545
$ws.Range("M26").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 876.8387
$ws.Range("I94").Value = 690.63635
$ws.Range("K94").Value = 690.63635
$ws.Range("M94").Value = -239.63635
$ws.Range("H118").Value = 7904.5454
$ws.Range("J118").Value = 7904.5454
$ws.Range("L118").Value = 7904.5454
$ws.Range("N118").Value = -11218.5454
$ws.Range("H134").Value = 3864.7222
$ws.Range("I134").Value = 2614.2258
$ws.Range("J134").Value = 5550.174
$ws.Range("K134").Value = 7842.6774
$ws.Range("L134").Value = 16650.522
$ws.Range("M134").Value = -5307.6774
$ws.Range("N134").Value = -21720.522

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 25536.363
$ws.Range("I93").Value = 12750
$ws.Range("J93").Value = 32842.855
$ws.Range("K93").Value = 12750
$ws.Range("L93").Value = 32842.855
$ws.Range("M93").Value = -10878
$ws.Range("N93").Value = -36586.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 32.666668
$ws.Range("I12").Value = 44.142857
$ws.Range("J12").Value = 26.928572
$ws.Range("K12").Value = 132.428571
$ws.Range("L12").Value = 80.78571599999999
$ws.Range("M12").Value = 40.57142899999999
$ws.Range("N12").Value = -426.785716
$ws.Range("H101").Value = 21000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 21000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 63000
$ws.Range("N101").Value = -67868
$ws.Range("H131").Value = 3513.8245
$ws.Range("I131").Value = 3713.7693
$ws.Range("J131").Value = 3454.75
$ws.Range("K131").Value = 11141.3079
$ws.Range("L131").Value = 10364.25
$ws.Range("M131").Value = -6101.3079
$ws.Range("N131").Value = -20444.25
$ws.Range("H134").Value = 5380.8696
$ws.Range("I134").Value = 2784
$ws.Range("K134").Value = 8352
$ws.Range("M134").Value = -3282
$ws.Range("M101").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3105.2334
$ws.Range("I126").Value = 2288.889
$ws.Range("J126").Value = 3455.0952
$ws.Range("K126").Value = 6866.667
$ws.Range("L126").Value = 10365.2856
$ws.Range("M126").Value = -4396.667
$ws.Range("N126").Value = -15305.2856
$ws.Range("H132").Value = 6277.8125
$ws.Range("I132").Value = 7322.609
$ws.Range("J132").Value = 3607.7778
$ws.Range("K132").Value = 21967.827
$ws.Range("L132").Value = 10823.3334
$ws.Range("M132").Value = -19437.827
$ws.Range("N132").Value = -15883.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10033.111
$ws.Range("I7").Value = 7462.25
$ws.Range("K7").Value = 7462.25
$ws.Range("M7").Value = -7350.25
$ws.Range("H40").Value = 3492.25
$ws.Range("I40").Value = 4586.7144
$ws.Range("J40").Value = 1960
$ws.Range("K40").Value = 4586.7144
$ws.Range("L40").Value = 1960
$ws.Range("M40").Value = -4450.7144
$ws.Range("N40").Value = -2232
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H122").Value = 12759.866
$ws.Range("I122").Value = 9914.286
$ws.Range("J122").Value = 15249.75
$ws.Range("K122").Value = 29742.858
$ws.Range("L122").Value = 45749.25
$ws.Range("M122").Value = -27292.858
$ws.Range("N122").Value = -50649.25
$ws.Range("H126").Value = 10033.111
$ws.Range("I126").Value = 7462.25
$ws.Range("K126").Value = 22386.75
$ws.Range("M126").Value = -19916.75
$ws.Range("H136").Value = 7113
$ws.Range("I136").Value = 2842.9524
$ws.Range("J136").Value = 37003.332
$ws.Range("K136").Value = 8528.8572
$ws.Range("L136").Value = 111009.996
$ws.Range("M136").Value = -5978.8572
$ws.Range("N136").Value = -116109.996
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5144.6313
$ws.Range("I62").Value = 5312.5
$ws.Range("J62").Value = 5022.5454
$ws.Range("K62").Value = 5312.5
$ws.Range("L62").Value = 5022.5454
$ws.Range("M62").Value = -4688.5
$ws.Range("N62").Value = -6270.5454
$ws.Range("H64").Value = 17208.348
$ws.Range("J64").Value = 17208.348
$ws.Range("L64").Value = 17208.348
$ws.Range("N64").Value = -17704.348
$ws.Range("H65").Value = 5144.6313
$ws.Range("I65").Value = 5312.5
$ws.Range("J65").Value = 5022.5454
$ws.Range("K65").Value = 26562.5
$ws.Range("L65").Value = 25112.727
$ws.Range("M65").Value = -23442.5
$ws.Range("N65").Value = -31352.727
$ws.Range("H67").Value = 17208.348
$ws.Range("J67").Value = 17208.348
$ws.Range("L67").Value = 17208.348
$ws.Range("N67").Value = -18924.348

Write-Host "Applied all Ultima Profits updates"